$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Joby Aviation)
$ws.Range("D2").Value = 15.76
$ws.Range("E2").Value = 57.5
$ws.Range("F2").Value = 11.61
$ws.Range("K2").Value = 57.9
$ws.Range("N2").Value = 54.85170003294819

# Row 3 (Archer Aviation)
$ws.Range("D3").Value = 8.949999999999999
$ws.Range("E3").Value = 61.4
$ws.Range("F3").Value = 19.49
$ws.Range("K3").Value = 54.5
$ws.Range("N3").Value = 54.85170003294819
